$wb = $excel.ActiveWorkbook

# =========================================================================
# "Generate Report for Handoff"
#
# The localization file previously tracked as
#   d511b3b9-99f0-42b3-ae4c-875eef848595.md
# was regenerated under a new id
#   a1b40c7b-e739-4b80-ad70-a9a93f3f31b8.md
# with a fresh handoff (new .xlf hashes + handoff timestamps), and the
# stale "Handoff transform failed" entry for
#   1d2fd01f-6e75-405c-a4e8-0a961974efe1.md
# is no longer part of the report (that row is removed on every sheet).
# =========================================================================

$oldMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/5f38f8aa8ae1c03bfbb3c6dd837be70df5f1b350/e2e/d511b3b9-99f0-42b3-ae4c-875eef848595.md"
$failedMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5f38f8aa8ae1c03bfbb3c6dd837be70df5f1b350/e2e/1d2fd01f-6e75-405c-a4e8-0a961974efe1.md"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ac99ce8643c0f000c1d773a8c343926f0e80fc7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/d511b3b9-99f0-42b3-ae4c-875eef848595.312d8a6dc63e6572b389ef896927e10464cd8f9e.zh-cn.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ec5bcc1630d8f26b9ea304d5a5b160541fcdc80/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/d511b3b9-99f0-42b3-ae4c-875eef848595.312d8a6dc63e6572b389ef896927e10464cd8f9e.de-de.xlf"

$newMdName = "a1b40c7b-e739-4b80-ad70-a9a93f3f31b8.md"
$newZhXlfName = "a1b40c7b-e739-4b80-ad70-a9a93f3f31b8.8fec1b573004ba8838c9bb103b3808b9c70d7f39.zh-cn.xlf"
$newDeXlfName = "a1b40c7b-e739-4b80-ad70-a9a93f3f31b8.8fec1b573004ba8838c9bb103b3808b9c70d7f39.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $newMdName

# Drop the "Handoff transform failed" row; the ".localization-config"
# row below it shifts up into row 3.
$wsOverview.Rows.Item(3).Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $oldMdUrl, "", "", $newMdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $failedMdUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A1").Hyperlinks.Delete()

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("C2").Value = $newZhXlfName
$wsZh.Range("D2").Value = "2016-01-27 02:54:57"

$wsZh.Rows.Item(3).Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $oldMdUrl, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", $newZhXlfName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $failedMdUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A1").Hyperlinks.Delete()

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("C2").Value = $newDeXlfName
$wsDe.Range("D2").Value = "2016-01-27 02:55:11"

$wsDe.Rows.Item(3).Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $oldMdUrl, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", $newDeXlfName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $failedMdUrl, "", "", ".localization-config") | Out-Null
